$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Stamp the destination L:O cells with the correct number formats before
#    we touch (and eventually clear) the source cells they are modeled on.
#    (PasteSpecial xlPasteFormats = -4122)
# ---------------------------------------------------------------------------

# Summary block (L5:O5 -> L2:O2)
$ws.Range("L5:O5").Copy()
$ws.Range("L2").PasteSpecial(-4122)

# "entrega Efvo" data block (L7:N12 -> L4:N9)
$ws.Range("L7:N7").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L8:N8").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("L9:N9").Copy()
$ws.Range("L6").PasteSpecial(-4122)
$ws.Range("L10:N10").Copy()
$ws.Range("L7").PasteSpecial(-4122)
$ws.Range("L11:N11").Copy()
$ws.Range("L8").PasteSpecial(-4122)
$ws.Range("L12:N12").Copy()
$ws.Range("L9").PasteSpecial(-4122)

# New deposit cells D/E (and F/G/H/I for row 57) - model their format on the
# already-populated D2:E2 pair.
$ws.Range("D2:E2").Copy()
$ws.Range("D57").PasteSpecial(-4122)
$ws.Range("F57").PasteSpecial(-4122)
$ws.Range("H57").PasteSpecial(-4122)
$ws.Range("D91").PasteSpecial(-4122)
$ws.Range("D135").PasteSpecial(-4122)
$ws.Range("D219").PasteSpecial(-4122)
$ws.Range("D224").PasteSpecial(-4122)
$ws.Range("D238").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Clear the old L:O cells that are being vacated (values+formats).
# ---------------------------------------------------------------------------
$ws.Range("L5:O5").Clear()
$ws.Range("L10:N12").Clear()

# ---------------------------------------------------------------------------
# 3) Write the relocated summary block at L2:O2.
# ---------------------------------------------------------------------------
$ws.Range("L2").Value = 341
$ws.Range("M2").Formula = '=SUMIF(D:I,"<>0")'
$ws.Range("N2").Value = 9000000
$ws.Range("O2").Formula = "=(N2-M2)*-1"

# ---------------------------------------------------------------------------
# 4) Write the relocated "entrega Efvo" data block at L4:N9.
# ---------------------------------------------------------------------------
$ws.Range("L4").Value = "1er entrega Efvo"
$ws.Range("M4").Value = 355000
$ws.Range("N4").Value = "18/01/2026"

$ws.Range("L5").Value = "2da entrega Efvo"
$ws.Range("M5").Value = 230000
$ws.Range("N5").Value = "20/01/2026"

$ws.Range("L6").Value = "3ra entrega Efvo"
$ws.Range("M6").Value = 230000
$ws.Range("N6").Value = "22/01/2026"

$ws.Range("L7").Value = "4ta entrega Efvo"
$ws.Range("M7").Value = 150000
$ws.Range("N7").Value = "24/01/2026"

$ws.Range("L8").Value = "5ta entrega Efvo"
$ws.Range("M8").Value = 156000
$ws.Range("N8").Value = "26/01/2026"

$ws.Range("L9").Value = "6ta entrega Efvo"
$ws.Range("M9").Value = 110000
$ws.Range("N9").Value = "26/01/2026"

# ---------------------------------------------------------------------------
# 5) New deposit values for rows 57, 91, 135, 219, 224, 238 (Transferencia).
# ---------------------------------------------------------------------------
$ws.Range("D57").Value = 10000
$ws.Range("E57").Value = "Transferencia"
$ws.Range("F57").Value = 10000
$ws.Range("G57").Value = "Transferencia"
$ws.Range("H57").Value = 10000
$ws.Range("I57").Value = "Transferencia"
$ws.Range("J57").Formula = "=SUM(D57:I57)"

$ws.Range("D91").Value = 10000
$ws.Range("E91").Value = "Transferencia"

$ws.Range("D135").Value = 10000
$ws.Range("E135").Value = "Transferencia"

$ws.Range("D219").Value = 10000
$ws.Range("E219").Value = "Transferencia"

$ws.Range("D224").Value = 10000
$ws.Range("E224").Value = "Transferencia"

$ws.Range("D238").Value = 10000
$ws.Range("E238").Value = "Transferencia"

# Row 248: existing deposit corrected from 5000 to 10000.
$ws.Range("D248").Value = 10000

# ---------------------------------------------------------------------------
# 6) Conditional formatting: the L2:M2 rules (3 of them) are removed as part
#    of this edit, and the remaining two rules' priorities shift down.
# ---------------------------------------------------------------------------
$cfCount = $ws.Cells.FormatConditions.Count
for ($i = $cfCount; $i -ge 1; $i--) {
    $fc = $ws.Cells.FormatConditions.Item($i)
}
